# Apply the "Updating Baselines and Daily Scores" edit to Daily_Scores
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update existing cells with newly recomputed scores (rows 3-37) ---
# Row 3
$ws.Range("E3").Value = 7.032520325203253
$ws.Range("G3").Value = 7.842595042636203
$ws.Range("J3").Value = 7.784216772999354
$ws.Range("K3").Value = 7.690726035167328
$ws.Range("M3").Value = 6.736111111111111
$ws.Range("P3").Value = 54.3019525141179
$ws.Range("Q3").Value = 47.78421677299936

# Row 5
$ws.Range("E5").Value = 0
$ws.Range("J5").Value = 9.665356265356268
$ws.Range("L5").Value = 8.88682541623718
$ws.Range("M5").Value = 10
$ws.Range("P5").Value = 10
$ws.Range("Q5").Value = 38.55218168159345

# Row 7
$ws.Range("F7").Value = 5
$ws.Range("K7").Value = 7.604166666666667
$ws.Range("M7").Value = 5.475474389477678
$ws.Range("P7").Value = 48.07964105614434
$ws.Range("Q7").Value = 35

# Row 9
$ws.Range("D9").Value = 7.299645419586821
$ws.Range("G9").Value = 7.595868644067797
$ws.Range("H9").Value = 9.02139307672871
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 9.054870827848864
$ws.Range("L9").Value = 8.624766330648681
$ws.Range("M9").Value = 0
$ws.Range("P9").Value = 26.65073947191666
$ws.Range("Q9").Value = 24.94580482696421

# Row 11
$ws.Range("C11").Value = 5.326797385620915
$ws.Range("F11").Value = 6.08702647196255
$ws.Range("G11").Value = 0
$ws.Range("I11").Value = 6.419094569740111
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 7.599735299678578
$ws.Range("M11").Value = 7.294146825396826
$ws.Range("P11").Value = 31.63977408043643
$ws.Range("Q11").Value = 26.08702647196255

# Row 13
$ws.Range("C13").Value = 9.947515212981751
$ws.Range("F13").Value = 7.501687341759451
$ws.Range("J13").Value = 7.684029484029486
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 8.130284078064943
$ws.Range("P13").Value = 28.07779929104669
$ws.Range("Q13").Value = 15.18571682578894

# Row 15
$ws.Range("C15").Value = 0
$ws.Range("E15").Value = 8.013737033922062
$ws.Range("F15").Value = 8.378595791234956
$ws.Range("G15").Value = 6.461988304093567
$ws.Range("I15").Value = 5.744047619047619
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").Value = 6.862649821833496
$ws.Range("P15").Value = 32.08242277889674
$ws.Range("Q15").Value = 38.37859579123496

# Row 17
$ws.Range("C17").Value = 8.56883874239351
$ws.Range("D17").Value = 7.09988180652894
$ws.Range("F17").Value = 8.906411898685914
$ws.Range("G17").Value = 7.744835805084748
$ws.Range("H17").Value = 8.371659587780195
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 9.995655062285692
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 0
$ws.Range("P17").Value = 26.30932960976395
$ws.Range("Q17").Value = 24.37795329299505

# Row 19
$ws.Range("C19").Value = 7.577812075111158
$ws.Range("E19").Value = 5.560695262125035
$ws.Range("J19").Value = 6.75
$ws.Range("M19").Value = 0
$ws.Range("P19").Value = 28.13850733723619
$ws.Range("Q19").Value = 46.75

# Row 21
$ws.Range("C21").Value = 7.142621703853955
$ws.Range("E21").Value = 7
$ws.Range("J21").Value = 7.259459459459459
$ws.Range("L21").Value = 7.524118171176993
$ws.Range("M21").Value = 10
$ws.Range("P21").Value = 34.14262170385395
$ws.Range("Q21").Value = 14.78357763063645

# Row 23
$ws.Range("C23").Value = 7.620071289457242
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 8.506284634088708
$ws.Range("J23").Value = 7.583333333333333
$ws.Range("K23").Value = 7.060290224995273
$ws.Range("M23").Value = 5.242954324586978
$ws.Range("P23").Value = 33.4296004731282
$ws.Range("Q23").Value = 37.58333333333333

# Row 25
$ws.Range("C25").Value = 8.331135902636916
$ws.Range("E25").Value = 8.1904998115042
$ws.Range("J25").Value = 8.533169533169534
$ws.Range("K25").Value = 8.163601553119236
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 10
$ws.Range("P25").Value = 34.68523726726035
$ws.Range("Q25").Value = 18.53316953316953

# Row 26
$ws.Range("I26").Value = 9.801639749334793
$ws.Range("P26").Value = 57.23294243563787

# Row 27
$ws.Range("E27").Value = 9.399880352993252
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 8.593941749940926
$ws.Range("J27").Value = 6
$ws.Range("K27").Value = 7.989695594630365
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 5.186793510151945
$ws.Range("P27").Value = 46.17031120771648
$ws.Range("Q27").Value = 31

# Row 29
$ws.Range("D29").Value = 9.197399743636698
$ws.Range("E29").Value = 7.773824877477729
$ws.Range("F29").Value = 7.020067493670377
$ws.Range("G29").Value = 7.893802966101695
$ws.Range("I29").Value = 10
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 9.153900747263267
$ws.Range("L29").Value = 7.209647268470797
$ws.Range("M29").Value = 0
$ws.Range("P29").Value = 34.82152859084269
$ws.Range("Q29").Value = 23.42711450577787

# Row 31
$ws.Range("E31").Value = 8.784693019343987
$ws.Range("G31").Value = 7.479444134105586
$ws.Range("I31").Value = 9.030257936507937
$ws.Range("J31").Value = 5.362488056318477
$ws.Range("P31").Value = 60.29439508995751
$ws.Range("Q31").Value = 45.36248805631848

# Row 33
$ws.Range("E33").Value = 7.416674934026471
$ws.Range("H33").Value = 7.866311318598018
$ws.Range("J33").Value = 10
$ws.Range("K33").Value = 7.272332278389609
$ws.Range("L33").Value = 9.358531770296477
$ws.Range("M33").Value = 10
$ws.Range("P33").Value = 34.68900721241608
$ws.Range("Q33").Value = 37.2248430888945

# Row 34
$ws.Range("C34").Value = 9.543878415821148
$ws.Range("E34").Value = 7.52149729394128
$ws.Range("F34").Value = 9.989386057808135
$ws.Range("M34").Value = 9.204521066022433
$ws.Range("P34").Value = 46.23174685189308
$ws.Range("Q34").Value = 39.98938605780813

# Row 35
$ws.Range("C35").Value = 8.512797192518279
$ws.Range("F35").Value = 7.298914092731612
$ws.Range("G35").Value = 8.143133725558091
$ws.Range("I35").Value = 6.557542820446463
$ws.Range("L35").Value = 10
$ws.Range("M35").Value = 5.713211584216516
$ws.Range("P35").Value = 33.92668532273935
$ws.Range("Q35").Value = 47.29891409273161

# Row 37
$ws.Range("C37").Value = 8.188514198782963
$ws.Range("G37").Value = 8.390360169491528
$ws.Range("J37").Value = 10
$ws.Range("K37").Value = 9.649050344335279
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = 10
$ws.Range("P37").Value = 56.22792471260977
$ws.Range("Q37").Value = 30

# --- Append new rows 38-41 for date 2025-02-10 ---
# Row 38
$ws.Range("A38").Value = "'2025-02-10"
$ws.Range("A38").Style = "Normal"
$ws.Range("B38").Value = "abs_activity"
$ws.Range("C38").Value = 8.827482417002898
$ws.Range("D38").Value = 0
$ws.Range("E38").Value = 10
$ws.Range("F38").Value = 4.408595154824254
$ws.Range("G38").Value = 8.108288623909186
$ws.Range("H38").Value = 10
$ws.Range("I38").Value = 9.828778776927553
$ws.Range("J38").Value = 4.851680761956861
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 7.92984742569226
$ws.Range("M38").Value = 8.320480120057807
$ws.Range("N38").Value = 0
$ws.Range("O38").Value = 0
$ws.Range("P38").Value = 45.08502993789745
$ws.Range("Q38").Value = 27.19012334247337

# Row 39
$ws.Range("A39").Value = "'2025-02-10"
$ws.Range("A39").Style = "Normal"
$ws.Range("B39").Value = "rel_activity"
$ws.Range("C39").Value = 6.875252636607502
$ws.Range("D39").Value = 5
$ws.Range("E39").Value = 10
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 10
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = 5.248015873015873
$ws.Range("N39").Value = 5
$ws.Range("O39").Value = 5
$ws.Range("P39").Value = 27.12326850962338
$ws.Range("Q39").Value = 20

# Row 40
$ws.Range("A40").Value = "'2025-02-10"
$ws.Range("A40").Style = "Normal"
$ws.Range("B40").Value = "abs_sleep"
$ws.Range("C40").Value = 10
$ws.Range("D40").Value = 0
$ws.Range("E40").Value = 10
$ws.Range("F40").Value = 8.566666666666666
$ws.Range("G40").Value = 10
$ws.Range("H40").Value = 10
$ws.Range("I40").Value = 10
$ws.Range("J40").Value = 7.199999999999999
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 3.733333333333334
$ws.Range("M40").Value = 9.133333333333333
$ws.Range("N40").Value = 0
$ws.Range("O40").Value = 0
$ws.Range("P40").Value = 49.13333333333333
$ws.Range("Q40").Value = 29.5

# Row 41
$ws.Range("A41").Value = "'2025-02-10"
$ws.Range("A41").Style = "Normal"
$ws.Range("B41").Value = "rel_sleep"
$ws.Range("C41").Value = 9.139325557809331
$ws.Range("D41").Value = 0
$ws.Range("E41").Value = 8.428599773805038
$ws.Range("F41").Value = 7.541822329100206
$ws.Range("G41").Value = 8.737950211864405
$ws.Range("H41").Value = 10
$ws.Range("I41").Value = 7.395212285456189
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = 0
$ws.Range("N41").Value = 0
$ws.Range("O41").Value = 0
$ws.Range("P41").Value = 33.70108782893496
$ws.Range("Q41").Value = 17.54182232910021

# --- Update sheet dimension to reflect the extended data range ---
# (Excel keeps this in sync automatically, nothing further required.)